$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing cells (value/type fixes) ---
$ws.Range("DT9").Value = 1
$ws.Range("BS11").Value = 1
$ws.Range("BS14").Value = "NaN"
$ws.Range("CU17").Value = "NaN"
$ws.Range("BQ19").Value = "NaN"
$ws.Range("CR19").Value = 1
$ws.Range("CC42").Value = "NaN"
$ws.Range("AP85").Value = "NaN"
$ws.Range("CF96").Value = "NaN"
$ws.Range("AI97").Value = "NaN"
$ws.Range("AI98").Value = "NaN"
$ws.Range("CF100").Value = "NaN"
$ws.Range("AI104").Value = "NaN"
$ws.Range("CR105").Value = 39
$ws.Range("CR106").Value = 54
$ws.Range("CR107").Value = 55
$ws.Range("CR108").Value = 56
$ws.Range("CR109").Value = 57
$ws.Range("CR110").Value = 58
$ws.Range("CR111").Value = 67
$ws.Range("CR112").Value = 67
$ws.Range("CR113").Value = 72
$ws.Range("BY114").Value = 153
$ws.Range("CR114").Value = 74
$ws.Range("BY115").Value = 166
$ws.Range("CR115").Value = 87
$ws.Range("AI116").Value = "NaN"
$ws.Range("BY116").Value = 194
$ws.Range("CR116").Value = 99
$ws.Range("AI117").Value = "NaN"
$ws.Range("BY117").Value = 211
$ws.Range("CR117").Value = 106
$ws.Range("AI118").Value = "NaN"
$ws.Range("BY118").Value = 220
$ws.Range("CR118").Value = 107
$ws.Range("BY119").Value = 228
$ws.Range("CR119").Value = 112
$ws.Range("BY120").Value = 246
$ws.Range("CR120").Value = 121
$ws.Range("BY121").Value = 249
$ws.Range("CR121").Value = 169
$ws.Range("BY122").Value = 273
$ws.Range("CR122").Value = 189
$ws.Range("BY123").Value = 289
$ws.Range("CR123").Value = 217
$ws.Range("BY124").Value = 305
$ws.Range("CR124").Value = 276
$ws.Range("BY125").Value = 325
$ws.Range("CR125").Value = 284
$ws.Range("BY126").Value = 345
$ws.Range("CR126").Value = 325
$ws.Range("BY127").Value = 374
$ws.Range("CR127").Value = 373
$ws.Range("AF128").Value = 12
$ws.Range("BY128").Value = 385
$ws.Range("CR128").Value = 402
$ws.Range("BY129").Value = 397
$ws.Range("CR129").Value = 428
$ws.Range("BY130").Value = 431
$ws.Range("CR130").Value = 456
$ws.Range("BY131").Value = 446
$ws.Range("CR131").Value = 460
$ws.Range("BY132").Value = 605
$ws.Range("CR132").Value = 474
$ws.Range("AF133").Value = 31
$ws.Range("BY133").Value = 637
$ws.Range("CR133").Value = 520
$ws.Range("AF134").Value = 33
$ws.Range("BY134").Value = 750
$ws.Range("CR134").Value = 552
$ws.Range("BY135").Value = 905
$ws.Range("CR135").Value = 561
$ws.Range("BY136").Value = 971
$ws.Range("CR136").Value = 597
$ws.Range("BY137").Value = 1061
$ws.Range("CR137").Value = 688
$ws.Range("BY138").Value = 1122
$ws.Range("CR138").Value = 725
$ws.Range("BY139").Value = 1218
$ws.Range("CR139").Value = 827
$ws.Range("BY140").Value = 1272
$ws.Range("CR140").Value = 875
$ws.Range("BY141").Value = 1295
$ws.Range("CR141").Value = 914
$ws.Range("BY142").Value = 1343
$ws.Range("CR142").Value = 996
$ws.Range("BY143").Value = 1372
$ws.Range("CR143").Value = 1013
$ws.Range("BY144").Value = 1407
$ws.Range("CR144").Value = 1160
$ws.Range("BY145").Value = 1527
$ws.Range("CR145").Value = 1185
$ws.Range("BY146").Value = 1586
$ws.Range("CR146").Value = 1217
$ws.Range("AF147").Value = 72
$ws.Range("BY147").Value = 1639
$ws.Range("CR147").Value = 1244
$ws.Range("BY148").Value = 1786
$ws.Range("CR148").Value = 1325
$ws.Range("BY149").Value = 1870
$ws.Range("CR149").Value = 1426
$ws.Range("BY150").Value = 1994
$ws.Range("CR150").Value = 1682
$ws.Range("BY151").Value = 2112
$ws.Range("CR151").Value = 1814
$ws.Range("BY152").Value = 2295
$ws.Range("CR152").Value = 1916
$ws.Range("AF153").Value = "NaN"
$ws.Range("BY153").Value = 2365
$ws.Range("CR153").Value = 1957
$ws.Range("AF154").Value = "NaN"
$ws.Range("BY154").Value = 2417
$ws.Range("CR154").Value = 2053
$ws.Range("AF155").Value = 176
$ws.Range("BY155").Value = 2487
$ws.Range("CR155").Value = 2128
$ws.Range("BY156").Value = 2584
$ws.Range("CR156").Value = 2186
$ws.Range("BY157").Value = 2650
$ws.Range("CR157").Value = 2233
$ws.Range("BY158").Value = 2725
$ws.Range("CR158").Value = 2287
$ws.Range("BY159").Value = 2829
$ws.Range("CR159").Value = 2333
$ws.Range("BY160").Value = 2886
$ws.Range("CR160").Value = 2389
$ws.Range("BY161").Value = 2963
$ws.Range("CR161").Value = 2435
$ws.Range("BY162").Value = 3049
$ws.Range("CR162").Value = 2483
$ws.Range("BY163").Value = 3123
$ws.Range("CR163").Value = 2542
$ws.Range("BY164").Value = 3243
$ws.Range("CR164").Value = 2569
$ws.Range("BY165").Value = 3294
$ws.Range("CR165").Value = 2607
$ws.Range("BY166").Value = 3362
$ws.Range("CR166").Value = 2636
$ws.Range("BY167").Value = 3443
$ws.Range("CR167").Value = 2652
$ws.Range("BY168").Value = 3509
$ws.Range("CR168").Value = 2693
$ws.Range("BY169").Value = 3564
$ws.Range("CR169").Value = 2725
$ws.Range("BY170").Value = 3634
$ws.Range("CR170").Value = 2764
$ws.Range("BY171").Value = 3672
$ws.Range("CR171").Value = 2769
$ws.Range("BY172").Value = 3773
$ws.Range("CR172").Value = 2816
$ws.Range("BY173").Value = 3818
$ws.Range("CR173").Value = 2821
$ws.Range("BY174").Value = 3857
$ws.Range("CR174").Value = 2835
$ws.Range("BY175").Value = 3920
$ws.Range("CR175").Value = 2861
$ws.Range("BY176").Value = 4041
$ws.Range("CR176").Value = 2917
$ws.Range("BY177").Value = 4101
$ws.Range("CR177").Value = 2929
$ws.Range("BY178").Value = 4148
$ws.Range("CR178").Value = 2959
$ws.Range("BY179").Value = 4195
$ws.Range("CR179").Value = 2966
$ws.Range("BY180").Value = 4234
$ws.Range("CR180").Value = 2974
$ws.Range("BY181").Value = 4281
$ws.Range("CR181").Value = 2976
$ws.Range("BY182").Value = 4380
$ws.Range("CR182").Value = 3000
$ws.Range("BY183").Value = 4661
$ws.Range("CR183").Value = 3054
$ws.Range("BY184").Value = 4805
$ws.Range("CR184").Value = 3061
$ws.Range("BY185").Value = 4860
$ws.Range("CR185").Value = 3067
$ws.Range("BY186").Value = 4884
$ws.Range("CR186").Value = 3088
$ws.Range("BY187").Value = 4934
$ws.Range("CR187").Value = 3092
$ws.Range("BY188").Value = 4963
$ws.Range("CR188").Value = 3099

# --- Append new data row 191 (date 2020-09-11 / serial 44085) ---
$ws.Cells.Item(191, 1).Value = 44085
$ws.Cells.Item(191, 2).Value = 702088
$ws.Cells.Item(191, 3).Value = 2721
$ws.Cells.Item(191, 4).Value = 93307
$ws.Cells.Item(191, 5).Value = 65631
$ws.Cells.Item(191, 6).Value = 236313
$ws.Cells.Item(191, 7).Value = 26786
$ws.Cells.Item(191, 8).Value = 4823
$ws.Cells.Item(191, 9).Value = 3805
$ws.Cells.Item(191, 10).Value = 7239
$ws.Cells.Item(191, 11).Value = 7012
$ws.Cells.Item(191, 12).Value = 14983
$ws.Cells.Item(191, 13).Value = 3853
$ws.Cells.Item(191, 14).Value = 21892
$ws.Cells.Item(191, 15).Value = 27572
$ws.Cells.Item(191, 16).Value = 6438
$ws.Cells.Item(191, 17).Value = 7139
$ws.Cells.Item(191, 18).Value = 13607
$ws.Cells.Item(191, 19).Value = 11125
$ws.Cells.Item(191, 20).Value = 15871
$ws.Cells.Item(191, 21).Value = 13429
$ws.Cells.Item(191, 22).Value = 3334
$ws.Cells.Item(191, 23).Value = 1823
$ws.Cells.Item(191, 24).Value = 7829
$ws.Cells.Item(191, 25).Value = 23720
$ws.Cells.Item(191, 26).Value = 12985
$ws.Cells.Item(191, 27).Value = 9246
$ws.Cells.Item(191, 28).Value = 52559
$ws.Cells.Item(191, 29).Value = 1432
$ws.Cells.Item(191, 30).Value = 431
$ws.Cells.Item(191, 31).Value = 519
$ws.Cells.Item(191, 32).Value = 458
$ws.Cells.Item(191, 33).Value = 340
$ws.Cells.Item(191, 34).Value = 223
$ws.Cells.Item(191, 35).Value = 442
$ws.Cells.Item(191, 36).Value = 2000
$ws.Cells.Item(191, 37).Value = 4143
$ws.Cells.Item(191, 38).Value = 36983
$ws.Cells.Item(191, 39).Value = 8073
$ws.Cells.Item(191, 40).Value = 2486
$ws.Cells.Item(191, 41).Value = 41074
$ws.Cells.Item(191, 42).Value = 1033
$ws.Cells.Item(191, 43).Value = 21545
$ws.Cells.Item(191, 44).Value = 1482
$ws.Cells.Item(191, 45).Value = 9254
$ws.Cells.Item(191, 46).Value = 1595
$ws.Cells.Item(191, 47).Value = 1585
$ws.Cells.Item(191, 48).Value = 6085
$ws.Cells.Item(191, 49).Value = 1753
$ws.Cells.Item(191, 50).Value = 954
$ws.Cells.Item(191, 51).Value = 2484
$ws.Cells.Item(191, 52).Value = 2652
$ws.Cells.Item(191, 53).Value = 53991
$ws.Cells.Item(191, 54).Value = 13272
$ws.Cells.Item(191, 55).Value = 4160
$ws.Cells.Item(191, 56).Value = 8566
$ws.Cells.Item(191, 57).Value = 5207
$ws.Cells.Item(191, 58).Value = 280
$ws.Cells.Item(191, 59).Value = 1433
$ws.Cells.Item(191, 60).Value = 2647
$ws.Cells.Item(191, 61).Value = 735
$ws.Cells.Item(191, 62).Value = 2110
$ws.Cells.Item(191, 63).Value = 9094
$ws.Cells.Item(191, 64).Value = 9082
$ws.Cells.Item(191, 65).Value = 9624
$ws.Cells.Item(191, 66).Value = 14060
$ws.Cells.Item(191, 67).Value = 1921
$ws.Cells.Item(191, 68).Value = 858
$ws.Cells.Item(191, 69).Value = 11215
$ws.Cells.Item(191, 70).Value = 8879
$ws.Cells.Item(191, 71).Value = 10373
$ws.Cells.Item(191, 72).Value = 1969
$ws.Cells.Item(191, 73).Value = 1791
$ws.Cells.Item(191, 74).Value = 4264
$ws.Cells.Item(191, 75).Value = 4064
$ws.Cells.Item(191, 76).Value = 1311
$ws.Cells.Item(191, 77).Value = 5165
$ws.Cells.Item(191, 78).Value = 2908
$ws.Cells.Item(191, 79).Value = 1669
$ws.Cells.Item(191, 80).Value = 826
$ws.Cells.Item(191, 81).Value = 2531
$ws.Cells.Item(191, 82).Value = 2116
$ws.Cells.Item(191, 83).Value = 1601
$ws.Cells.Item(191, 84).Value = 1244
$ws.Cells.Item(191, 85).Value = 5867
$ws.Cells.Item(191, 86).Value = 1788
$ws.Cells.Item(191, 87).Value = 1304
$ws.Cells.Item(191, 88).Value = 1499
$ws.Cells.Item(191, 89).Value = 1880
$ws.Cells.Item(191, 90).Value = 1808
$ws.Cells.Item(191, 91).Value = 2165
$ws.Cells.Item(191, 92).Value = 1355
$ws.Cells.Item(191, 93).Value = 1163
$ws.Cells.Item(191, 94).Value = 1156
$ws.Cells.Item(191, 95).Value = 729
$ws.Cells.Item(191, 96).Value = 3173
$ws.Cells.Item(191, 97).Value = 1246
$ws.Cells.Item(191, 98).Value = 857
$ws.Cells.Item(191, 99).Value = 891
$ws.Cells.Item(191, 100).Value = 1602
$ws.Cells.Item(191, 101).Value = 1405
$ws.Cells.Item(191, 102).Value = 709
$ws.Cells.Item(191, 103).Value = 811
$ws.Cells.Item(191, 104).Value = 1102
$ws.Cells.Item(191, 105).Value = 1400
$ws.Cells.Item(191, 106).Value = 1207
$ws.Cells.Item(191, 107).Value = 1321
$ws.Cells.Item(191, 108).Value = 1038
$ws.Cells.Item(191, 109).Value = 325
$ws.Cells.Item(191, 110).Value = 347
$ws.Cells.Item(191, 111).Value = 755
$ws.Cells.Item(191, 112).Value = 680
$ws.Cells.Item(191, 113).Value = 447
$ws.Cells.Item(191, 114).Value = 535
$ws.Cells.Item(191, 115).Value = 360
$ws.Cells.Item(191, 116).Value = 644
$ws.Cells.Item(191, 117).Value = 727
$ws.Cells.Item(191, 118).Value = 518
$ws.Cells.Item(191, 119).Value = 484
$ws.Cells.Item(191, 120).Value = 372
$ws.Cells.Item(191, 121).Value = 518
$ws.Cells.Item(191, 122).Value = 127445
$ws.Cells.Item(191, 123).Value = 297939
$ws.Cells.Item(191, 124).Value = 13730
$ws.Cells.Item(191, 125).Value = 128174
$ws.Cells.Item(191, 126).Value = 79295
$ws.Cells.Item(191, 127).Value = 37198
$ws.Cells.Item(191, 128).Value = 10913

# --- Update selection to reflect the newly added row, like the source edit ---
$ws.Range("A191").Select()

